$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.3312325492543664,
    -0.0099999993697217349,
    -0.0089999993608156359,
    -0.011999999832617902,
    -0.005999999371945286,
    -0.0059999993573143229,
    -0.01999999924413487,
    -0.019999999239132649,
    -0.0059999993461037349,
    -0.0059999993427481968,
    -0.0044999993545928874,
    -0.0059999993423804909,
    -0.0059999993411397057,
    -0.011999999293447416,
    -0.0059999993420714048,
    -0.0059999993439063815,
    -0.0059999993484902703,
    0.029589319393379832,
    -0.0089999993764475761,
    -0.052061164673743221,
    -0.0089999993601281858,
    -0.0089999993595428762,
    -0.0089999993549891855,
    -0.041999999082397821,
    -0.041999999077262373,
    -0.0059999993571935306,
    -0.0059999993563097931,
    -0.0059999993520056805,
    -0.011999999301805175,
    -0.019999999237671151,
    -0.014999999278924037,
    0.045495068428927254,
    -0.020596381311274392
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
